# Update "Phan chia cong viec" worksheet per latest edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Names" column (E) assignments.
$ws.Range("E3").Value = "sangdang"
$ws.Range("E5").Value = "sangdang"
$ws.Range("E7").Value = "sangdang"
$ws.Range("E8").Value = "namdao"
$ws.Range("E10").Value = "namdao"
$ws.Range("E12").Value = "namdao"

# Update "Start Build" dates (column B) to 2023-08-16 (serial 45154) for rows 8-13.
$ws.Range("B8").Value = 45154
$ws.Range("B9").Value = 45154
$ws.Range("B10").Value = 45154
$ws.Range("B11").Value = 45154
$ws.Range("B12").Value = 45154
$ws.Range("B13").Value = 45154

# Update the active selection to reflect the latest cursor position.
$ws.Range("C16").Select()
